$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.676.13'
$ws.Range("E2").Value = '  +3.83%  '
$ws.Range("D3").Value = '3.300.14'
$ws.Range("E3").Value = '  +1.20%  '
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '576.09'
$ws.Range("E5").Value = '  +0.41%  '
$ws.Range("D6").Value = '177.93'
$ws.Range("E6").Value = '  -1.53%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '0.585'
$ws.Range("E8").Value = '  +3.68%  '
$ws.Range("D9").Value = '3.294.65'
$ws.Range("E9").Value = '  +1.26%  '
$ws.Range("D10").Value = '0.175'
$ws.Range("E10").Value = '  +1.46%  '
$ws.Range("D11").Value = '0.574'
$ws.Range("E11").Value = '  +2.23%  '
$ws.Range("D12").Value = '45.71'
$ws.Range("E12").Value = '  +0.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000270'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.82%  '
$ws.Range("D14").Value = '702.24'
$ws.Range("E14").Value = '  +14.53%  '
$ws.Range("D15").Value = '3.833.39'
$ws.Range("E15").Value = '  +1.46%  '
$ws.Range("D16").Value = '8.36'
$ws.Range("E16").Value = '  +0.86%  '
$ws.Range("D17").Value = '67.800.94'
$ws.Range("E17").Value = '  +4.03%  '
$ws.Range("E18").Value = '  +1.73%  '
$ws.Range("D19").Value = '3.311.53'
$ws.Range("E19").Value = '  +1.50%  '
$ws.Range("D20").Value = '17.44'
$ws.Range("E20").Value = '  +0.03%  '
$ws.Range("D21").Value = '10.79'
$ws.Range("E21").Value = '  +0.41%  '
$ws.Range("D22").Value = '0.893'
$ws.Range("E22").Value = '  +1.98%  '
$ws.Range("D23").Value = '16.83'
$ws.Range("E23").Value = '  -7.31%  '
$ws.Range("D24").Value = '5.17'
$ws.Range("E24").Value = '  +5.39%  '
$ws.Range("D25").Value = '98.78'
$ws.Range("E25").Value = '  +0.66%  '
$ws.Range("D26").Value = '3.93'
$ws.Range("E26").Value = '  +0.54%  '
$ws.Range("D27").Value = '2.73'
$ws.Range("E27").Value = '  +1.81%  '
$ws.Range("D28").Value = '9.34'
$ws.Range("E28").Value = '  +0.85%  '
$ws.Range("D29").Value = '33.11'
$ws.Range("E29").Value = '  +9.73%  '
$ws.Range("D30").Value = '8.46'
$ws.Range("E30").Value = '  +2.70%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.70'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '578.90'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.01%  '
$ws.Range("D33").Value = '3.912.61'
$ws.Range("E33").Value = '  +4.37%  '
$ws.Range("D34").Value = '10.84'
$ws.Range("E34").Value = '  +1.39%  '
$ws.Range("E35").Value = '  +1.84%  '
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  +0.19%  '
$ws.Range("D37").Value = '3.36'
$ws.Range("E37").Value = '  -5.36%  '
$ws.Range("D38").Value = '55.29'
$ws.Range("E38").Value = '  -0.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.130'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.41%  '
$ws.Range("D40").Value = '3.15'
$ws.Range("E40").Value = '  +2.06%  '
$ws.Range("D41").Value = '2.61'
$ws.Range("E41").Value = '  +2.91%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '32.10'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.15%  '
$ws.Range("D43").Value = '0.0₃0680'
$ws.Range("E43").Value = '  +1.94%  '
$ws.Range("D44").Value = '3.33'
$ws.Range("E44").Value = '  -1.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.330'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.67%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0410'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.73%  '
$ws.Range("D47").Value = '0.128'
$ws.Range("E47").Value = '  +2.90%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +12.17%  '
$ws.Range("E49").Value = '  +0.24%  '
$ws.Range("D50").Value = '2.52'
$ws.Range("E50").Value = '  +2.58%  '
$ws.Range("D51").Value = '128.25'
$ws.Range("E51").Value = '  +0.41%  '
